$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 924, pushing the existing rows 924:994 down to 926:996.
$ws.Rows("924:925").Insert()

# Row 924 - new record
$ws.Range("A924").Value = 5
$ws.Range("B924").Value = "Macroferia Regional de Talca"
$ws.Range("C924").Value = "Maule"
$ws.Range("D924").Value = 45265
$ws.Range("E924").Value = 7
$ws.Range("F924").Value = 100114001
$ws.Range("G924").Value = "Papa"
$ws.Range("H924").Value = "Asterix"
$ws.Range("I924").Value = "1a nueva(o)"
$ws.Range("J924").Value = 1500
$ws.Range("K924").Value = 20000
$ws.Range("L924").Value = 20000
$ws.Range("M924").Value = 20000
$ws.Range("N924").Value = "`$/saco 25 kilos"
$ws.Range("O924").Value = "Región de O'Higgins"
$ws.Range("P924").Value = 800
$ws.Range("Q924").Value = 25
$ws.Range("R924").Value = "Hortaliza"

# Row 925 - new record
$ws.Range("A925").Value = 5
$ws.Range("B925").Value = "Macroferia Regional de Talca"
$ws.Range("C925").Value = "Maule"
$ws.Range("D925").Value = 45265
$ws.Range("E925").Value = 7
$ws.Range("F925").Value = 100114001
$ws.Range("G925").Value = "Papa"
$ws.Range("H925").Value = "Asterix"
$ws.Range("I925").Value = "1a nueva(o)"
$ws.Range("J925").Value = 1600
$ws.Range("K925").Value = 18000
$ws.Range("L925").Value = 18000
$ws.Range("M925").Value = 18000
$ws.Range("N925").Value = "`$/saco 25 kilos"
$ws.Range("O925").Value = "Región del Maule"
$ws.Range("P925").Value = 720
$ws.Range("Q925").Value = 25
$ws.Range("R925").Value = "Hortaliza"
